# Update crypto price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.699.55"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "3.052.93"
$ws.Range("E3").Value = "  -5.13%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "614.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.363"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.801"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +15.10%  "
$ws.Range("D10").Value = "3.049.78"
$ws.Range("E10").Value = "  -5.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.589"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -10.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").Value = "87.694.47"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "3.617.58"
$ws.Range("E16").Value = "  -5.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.63%  "
$ws.Range("D18").Value = "3.055.30"
$ws.Range("E18").Value = "  -5.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000193"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -18.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "416.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.63%  "
$ws.Range("E23").Value = "  -7.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "80.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").Value = "3.217.28"
$ws.Range("E28").Value = "  -5.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.161"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  -7.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "496.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -16.98%  "
$ws.Range("E35").Value = "  -7.81%  "
$ws.Range("E36").Value = "  -7.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E43").Value = "  -5.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("E45").Value = "  -8.64%  "
$ws.Range("E46").Value = "  +5.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0660"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -11.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.692"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.72%  "
$ws.Range("E51").Value = "  -8.53%  "
